$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every
# Avverkningsanmälan row. It was bumped by one day (serial 45180 ->
# 45181, i.e. 2023-09-11 -> 2023-09-12) across all data rows (2..176).
$ws.Range("C2:C176").Value = 45181
